$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "90.132.59"
$ws.Range("E2").Value = "  -1.17%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.075.80"
$ws.Range("E3").Value = "  -1.85%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "

# Row 5 - Solana
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.98"
$ws.Range("E5").Value = "  -0.61%  "

# Row 6 - BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "612.61"
$ws.Range("E6").Value = "  -1.75%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.06"
$ws.Range("E7").Value = "  +15.77%  "

# Row 8 - Dogecoin
$ws.Range("E8").Value = "  -7.24%  "

# Row 9 - USDC
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.07%  "

# Row 10 - LidoStakedEther
$ws.Range("D10").Value = "3.073.80"
$ws.Range("E10").Value = "  -1.76%  "

# Row 11 - Cardano
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.715"
$ws.Range("E11").Value = "  -4.12%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +2.64%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -5.68%  "

# Row 14 - Toncoin (was Avalanche)
$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.50"
$ws.Range("E14").Value = "  +1.43%  "

# Row 15 - Avalanche (was Toncoin)
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.28"
$ws.Range("E15").Value = "  +0.37%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "89.943.61"
$ws.Range("E16").Value = "  -1.20%  "

# Row 17 - WrappedliquidstakedEther2.0
$ws.Range("D17").Value = "3.631.20"
$ws.Range("E17").Value = "  -2.10%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.066.13"
$ws.Range("E18").Value = "  -2.73%  "

# Row 19 - SuiNetwork
$ws.Range("E19").Value = "  -4.11%  "

# Row 20 - Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.15"
$ws.Range("E20").Value = "  -0.29%  "

# Row 21 - PEPE
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000203"
$ws.Range("E21").Value = "  -10.68%  "

# Row 22 - BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "445.35"

# Row 23 - Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.84"
$ws.Range("E23").Value = "  +0.65%  "

# Row 24 - Polkadot
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.38"
$ws.Range("E24").Value = "  +4.04%  "

# Row 25 - NEARProtocol
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.95"
$ws.Range("E25").Value = "  -2.62%  "

# Row 26 - Litecoin
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "90.58"
$ws.Range("E26").Value = "  +8.03%  "

# Row 27 - Aptos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.88"
$ws.Range("E27").Value = "  -4.42%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "3.232.29"
$ws.Range("E28").Value = "  -1.71%  "

# Row 29 - Dai
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.02%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.31"
$ws.Range("E30").Value = "  +2.76%  "

# Row 31 - Cronos
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.158"
$ws.Range("E31").Value = "  -5.38%  "

# Row 32 - Binance-PegBSC-USD
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  +1.14%  "

# Row 33 - EthereumClassic
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.84"
$ws.Range("E33").Value = "  +18.63%  "

# Row 34 - Stellar
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.202"
$ws.Range("E34").Value = "  +39.86%  "

# Row 35 - Kaspa
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.148"
$ws.Range("E35").Value = "  +5.09%  "

# Row 36 - Bittensor
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "496.95"
$ws.Range("E36").Value = "  -6.16%  "

# Row 37 - PancakeSwap
$ws.Range("E37").Value = "  +1.16%  "

# Row 38 - RenderToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.81"
$ws.Range("E38").Value = "  -7.79%  "

# Row 39 - Fetch.AI
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.29"
$ws.Range("E39").Value = "  -2.00%  "

# Row 40 - dogwifhat
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.42"
$ws.Range("E40").Value = "  -11.91%  "

# Row 41 - PolygonEcosystemToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.426"
$ws.Range("E41").Value = "  +12.50%  "

# Row 42 - WhiteBITCoin
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.19"
$ws.Range("E42").Value = "  -0.43%  "

# Row 43 - Hedera
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0858"
$ws.Range("E43").Value = "  +11.20%  "

# Row 44 - USDe
$ws.Range("E44").Value = "  +0.01%  "

# Row 45 - Stacks
$ws.Range("E45").Value = "  +0.59%  "

# Row 46 - MantraDAO
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.00"
$ws.Range("E46").Value = "  +17.78%  "

# Row 47 - ARBITRUM (was Monero)
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.690"
$ws.Range("E47").Value = "  +11.45%  "

# Row 48 - Filecoin (was ARBITRUM)
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.59"
$ws.Range("E48").Value = "  +9.89%  "

# Row 49 - Monero (was Filecoin)
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "148.09"
$ws.Range("E49").Value = "  +3.08%  "

# Row 50 - OKB
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.44"
$ws.Range("E50").Value = "  +0.57%  "

# Row 51 - ImmutableX
$ws.Range("E51").Value = "  +1.43%  "
